# edit.ps1 - applies the syllabus date/content refresh described by the diff
# (Fall 2021-2022 -> Fall 2023-2024 schedule update, contact info refresh, etc.)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 4: title card - semester label + drop the "Download WORD, PDF" line
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange
$tr4.Paragraphs(3).Text = "Fall Semester, 2023-2024"
$tr4.Paragraphs(5).Delete()

# ---------------------------------------------------------------------------
# Slide 5: instructor info table
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tbl5 = $s5.Shapes.Item(1).Table
$tbl5.Cell(4, 2).Shape.TextFrame.TextRange.Text = "TBD"

$newRow5 = $tbl5.Rows.Add(5)
$tbl5.Cell(5, 1).Shape.TextFrame.TextRange.Text = "Microsoft Teams Code"
$tbl5.Cell(5, 1).Shape.TextFrame.TextRange.Font.Bold = -1
$tbl5.Cell(5, 2).Shape.TextFrame.TextRange.Text = "etj1k7b"

$tbl5.Cell(6, 2).Shape.TextFrame.TextRange.Text = "Wednesday 09:00 - 12:00 (Theory) / Thursday 13:00 - 14:30 (Lab)"

# ---------------------------------------------------------------------------
# Slide 6: classroom + office hours table
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tbl6 = $s6.Shapes.Item(1).Table
$tbl6.Cell(1, 2).Shape.TextFrame.TextRange.Text = "İİBF-414 (Level-4)"
$tbl6.Cell(2, 2).Shape.TextFrame.TextRange.Text = "Scheduled through your university account, meetings will take place via Google Meet and are coordinated by request through email. To expedite responses, please initiate your email subject line with the [CE103] tag and maintain a formal, concise, and clear email body."

# ---------------------------------------------------------------------------
# Slide 9: course description paragraph
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(2)
$shp9.TextFrame.TextRange.Paragraphs(1).Text = "The objective of this course is to systematically build foundational skills in algorithms and programming, aimed at enhancing students" + [char]0x2019 + " career prospects. The instructional approach focuses on the transfer of expert knowledge while guiding students in identifying effective learning techniques and practical applications within the realm of algorithms and programming. Emphasis is placed on hands-on projects and applications, fortifying the learning experience through practice over pure theory. The course covers functional programming in C, C++, C#, and Java, utilizing the latest development environments."

# ---------------------------------------------------------------------------
# Slides 39-44: weekly schedule tables (dates 2021 -> 2023, some -> TBD)
# ---------------------------------------------------------------------------
$s39 = $p.Slides.Item(39)
$tbl39 = $s39.Shapes.Item(2).Table
$tbl39.Cell(2, 2).Shape.TextFrame.TextRange.Text = "04.10.2023 05.10.2023"
$tbl39.Cell(3, 2).Shape.TextFrame.TextRange.Text = "11.10.2023 12.10.2023"

$s40 = $p.Slides.Item(40)
$tbl40 = $s40.Shapes.Item(1).Table
$tbl40.Cell(1, 2).Shape.TextFrame.TextRange.Text = "18.10.2023 19.10.2023"
$tbl40.Cell(2, 2).Shape.TextFrame.TextRange.Text = "25.10.2023 26.10.2023"
$tbl40.Cell(3, 2).Shape.TextFrame.TextRange.Text = "01.11.2023 02.11.2023"

$s41 = $p.Slides.Item(41)
$tbl41 = $s41.Shapes.Item(1).Table
$tbl41.Cell(1, 2).Shape.TextFrame.TextRange.Text = "08.11.2023 09.11.2023"
$tbl41.Cell(2, 2).Shape.TextFrame.TextRange.Text = "15.11.2023 16.11.2023"
$tbl41.Cell(3, 2).Shape.TextFrame.TextRange.Text = "TBD"

$s42 = $p.Slides.Item(42)
$tbl42 = $s42.Shapes.Item(1).Table
$tbl42.Cell(1, 2).Shape.TextFrame.TextRange.Text = "29.11.2023 30.11.2023"
$tbl42.Cell(2, 2).Shape.TextFrame.TextRange.Text = "06.11.2023 07.11.2023"
$tbl42.Cell(3, 2).Shape.TextFrame.TextRange.Text = "13.11.2023 14.11.2023"

$s43 = $p.Slides.Item(43)
$tbl43 = $s43.Shapes.Item(1).Table
$tbl43.Cell(1, 2).Shape.TextFrame.TextRange.Text = "20.11.2023 21.11.2023"
$tbl43.Cell(2, 2).Shape.TextFrame.TextRange.Text = "27.11.2023 28.11.2023"
$tbl43.Cell(3, 2).Shape.TextFrame.TextRange.Text = "03.12.2023 04.12.2023"

$s44 = $p.Slides.Item(44)
$tbl44 = $s44.Shapes.Item(1).Table
$tbl44.Cell(1, 2).Shape.TextFrame.TextRange.Text = "10.12.2023 11.12.2023"
$tbl44.Cell(2, 2).Shape.TextFrame.TextRange.Text = "TBD"

# Final row: remove the stray "    ." run trailing the bold "Final" label.
$tbl44.Rows(2).Delete()
$tbl44.Rows.Add() | Out-Null
$tbl44.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Week-16"
$tbl44.Cell(2, 2).Shape.TextFrame.TextRange.Text = "TBD"
$tbl44.Cell(2, 3).Shape.TextFrame.TextRange.Text = "Final"
$tbl44.Cell(2, 3).Shape.TextFrame.TextRange.Font.Bold = -1
$tbl44.Cell(2, 4).Shape.TextFrame.TextRange.Text = "TBD"
